# The old question #4 ("4. 다음의 뜻 풀이 중 옳지 못한 것은?" with its unused
# "가가" answer choices) is replaced by what used to be question #11 at the
# bottom of the sheet ("좋은거 있으면 농갈라무야지"), renumbered to #4 since it
# now takes that question's place. The old row 12 (now duplicated into row 5)
# is then removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pull row 12's question/answers up into row 5, renumbering "11." -> "4.".
$ws.Range("B5").Value = $ws.Range("B12").Value2
$ws.Range("C5").Value = $ws.Range("C12").Value2
$ws.Range("D5").Value = $ws.Range("D12").Value2
$ws.Range("E5").Value = $ws.Range("E12").Value2
$ws.Range("F5").Value = $ws.Range("F12").Value2
$ws.Range("A5").Value = "4. 다음 문장의 뜻을 올바르게 해석한 것은?`r`n좋은거 있으면 농갈라무야지"

# Row 5 now holds a question, so give it the same wrapped/taller style as the
# other question rows (matches cells like A3, A4, A7 ...).
$ws.Range("A5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 33

# The content that used to live in row 12 now lives in row 5, so drop row 12.
$ws.Rows.Item(12).Delete()

# Reflect the author's last selection in the saved view.
$ws.Range("D13").Select() | Out-Null
